# Refresh the Market Board price snapshot + recomputed Leve profit figures
# across every job sheet (the scheduled pricing-data runner pass).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1831.0834
$ws.Range("I2").Value = 1571.75
$ws.Range("J2").Value = 2349.75
$ws.Range("K2").Value = 1571.75
$ws.Range("L2").Value = 2349.75
$ws.Range("M2").Value = -1458.75
$ws.Range("N2").Value = -2575.75
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 2473.4
$ws.Range("J17").Value = 2624.4285
$ws.Range("L17").Value = 7873.2855
$ws.Range("N17").Value = -8209.2855
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 84289
$ws.Range("I107").Value = 143522.86
$ws.Range("K107").Value = 143522.86
$ws.Range("M107").Value = -141602.86
# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 304637.62
$ws.Range("I112").Value = 1750
$ws.Range("J112").Value = 324178.78
$ws.Range("K112").Value = 5250
$ws.Range("L112").Value = 972536.3400000001
$ws.Range("M112").Value = -4142
$ws.Range("N112").Value = -974752.3400000001
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2121.4285
$ws.Range("J138").Value = 2696.9556
$ws.Range("L138").Value = 8090.8668
$ws.Range("N138").Value = -18370.8668

$ws = $wb.Worksheets.Item("ARM")
# Row 43 (Leve Item ID 21715)
$ws.Range("H43").Value = 25372
$ws.Range("J43").Value = 25372
$ws.Range("L43").Value = 25372
$ws.Range("N43").Value = -25998
# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 55 (Leve Item ID 2830)
$ws.Range("H55").Value = 32500
$ws.Range("I55").Value = 7000
$ws.Range("K55").Value = 7000
$ws.Range("M55").Value = -6685
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3767.5789
$ws.Range("I74").Value = 2207.5
$ws.Range("J74").Value = 5501
$ws.Range("K74").Value = 2207.5
$ws.Range("L74").Value = 5501
$ws.Range("M74").Value = -1333.5
$ws.Range("N74").Value = -7249
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3767.5789
$ws.Range("I77").Value = 2207.5
$ws.Range("J77").Value = 5501
$ws.Range("K77").Value = 11037.5
$ws.Range("L77").Value = 27505
$ws.Range("M77").Value = -6669.5
$ws.Range("N77").Value = -36241
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 14230.467
$ws.Range("I102").Value = 10599.8
$ws.Range("K102").Value = 10599.8
$ws.Range("M102").Value = -8977.799999999999
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3808.7222
$ws.Range("I132").Value = 3403.3333
$ws.Range("K132").Value = 10209.9999
$ws.Range("M132").Value = -7679.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 4240
$ws.Range("I105").Value = 3500
$ws.Range("J105").Value = 4980
$ws.Range("K105").Value = 3500
$ws.Range("L105").Value = 4980
$ws.Range("M105").Value = -1753
$ws.Range("N105").Value = -8474

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1590.3
$ws.Range("I132").Value = 1326.1143
$ws.Range("K132").Value = 3978.3429
$ws.Range("M132").Value = -1448.3429

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 2512.3125
$ws.Range("I5").Value = 1022.1111
$ws.Range("K5").Value = 3066.3333
$ws.Range("M5").Value = -2954.3333
# Row 9 (Leve Item ID 4681)
$ws.Range("H9").Value = 3333336
$ws.Range("I9").Value = 5000001.5
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 15000004.5
$ws.Range("L9").Value = 15
$ws.Range("M9").Value = -14999780.5
$ws.Range("N9").Value = -463
# Row 10 (Leve Item ID 4689)
$ws.Range("H10").Value = 353.25
$ws.Range("I10").Value = 400.85715
$ws.Range("K10").Value = 1202.57145
$ws.Range("M10").Value = -1063.57145
# Row 20 (Leve Item ID 4651)
$ws.Range("H20").Value = 250
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 300
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -73
$ws.Range("N20").Value = -1654
# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 407
$ws.Range("I46").Value = 167.66667
$ws.Range("J46").Value = 1125
$ws.Range("K46").Value = 503.00001
$ws.Range("L46").Value = 3375
$ws.Range("M46").Value = -412.00001
$ws.Range("N46").Value = -3557
# Row 51 (Leve Item ID 4646)
$ws.Range("H51").Value = 14117.75
$ws.Range("I51").Value = 1997.5
$ws.Range("J51").Value = 26238
$ws.Range("K51").Value = 5992.5
$ws.Range("L51").Value = 78714
$ws.Range("M51").Value = -5532.5
$ws.Range("N51").Value = -79634
# Row 58 (Leve Item ID 4703)
$ws.Range("H58").Value = 40001
$ws.Range("J58").Value = 40001
$ws.Range("L58").Value = 120003
$ws.Range("N58").Value = -120259
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 2184.75
$ws.Range("I68").Value = 2613.3333
$ws.Range("J68").Value = 899
$ws.Range("K68").Value = 7839.999899999999
$ws.Range("L68").Value = 2697
$ws.Range("M68").Value = -7028.999899999999
$ws.Range("N68").Value = -4319
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 2184.75
$ws.Range("I71").Value = 2613.3333
$ws.Range("J71").Value = 899
$ws.Range("K71").Value = 23519.9997
$ws.Range("L71").Value = 8091
$ws.Range("M71").Value = -19463.9997
$ws.Range("N71").Value = -16203
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 302.08334
$ws.Range("I107").Value = 184.4
$ws.Range("J107").Value = 386.14285
$ws.Range("K107").Value = 553.2
$ws.Range("L107").Value = 1158.42855
$ws.Range("M107").Value = 1366.8
$ws.Range("N107").Value = -4998.428550000001
# Row 112 (Leve Item ID 27855)
$ws.Range("H112").Value = 5599.95
$ws.Range("J112").Value = 5599.95
$ws.Range("L112").Value = 16799.85
$ws.Range("N112").Value = -19015.85
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 2087.375
$ws.Range("I122").Value = 971.2857
$ws.Range("J122").Value = 9900
$ws.Range("K122").Value = 8741.5713
$ws.Range("L122").Value = 89100
$ws.Range("M122").Value = -6291.5713
$ws.Range("N122").Value = -94000
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 1680.7273
$ws.Range("I129").Value = 946.1111
$ws.Range("J129").Value = 4986.5
$ws.Range("K129").Value = 2838.3333
$ws.Range("L129").Value = 14959.5
$ws.Range("M129").Value = 2161.6667
$ws.Range("N129").Value = -24959.5
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 2512.3125
$ws.Range("I135").Value = 1022.1111
$ws.Range("K135").Value = 9198.999899999999
$ws.Range("M135").Value = -6663.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 260897.25
$ws.Range("I70").Value = 341863
$ws.Range("J70").Value = 18000
$ws.Range("K70").Value = 341863
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -341593
$ws.Range("N70").Value = -18540
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 260897.25
$ws.Range("I73").Value = 341863
$ws.Range("J73").Value = 18000
$ws.Range("K73").Value = 341863
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -340927
$ws.Range("N73").Value = -19872
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 716.7273
$ws.Range("I107").Value = 650
$ws.Range("K107").Value = 650
$ws.Range("M107").Value = 1270
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2964.2856
$ws.Range("I132").Value = 2250
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2322.8333
$ws.Range("I22").Value = 1131
$ws.Range("K22").Value = 1131
$ws.Range("M22").Value = -836
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2322.8333
$ws.Range("I27").Value = 1131
$ws.Range("K27").Value = 1131
$ws.Range("M27").Value = -1024
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1575.9546
$ws.Range("I55").Value = 1060.48
$ws.Range("K55").Value = 1060.48
$ws.Range("M55").Value = -887.48
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3125
$ws.Range("I68").Value = 3125
$ws.Range("K68").Value = 3125
$ws.Range("M68").Value = -2376
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3125
$ws.Range("I71").Value = 3125
$ws.Range("K71").Value = 15625
$ws.Range("M71").Value = -11881
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2517.35
$ws.Range("I93").Value = 2145.7273
$ws.Range("J93").Value = 2971.5557
$ws.Range("K93").Value = 2145.7273
$ws.Range("L93").Value = 2971.5557
$ws.Range("M93").Value = -897.7273
$ws.Range("N93").Value = -5467.5557
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2961.4075
$ws.Range("I132").Value = 2484.375
$ws.Range("K132").Value = 7453.125
$ws.Range("M132").Value = -4923.125
# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 6505.722
$ws.Range("I136").Value = 9162
$ws.Range("K136").Value = 27486
$ws.Range("M136").Value = -24936

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (Leve Item ID 2670)
$ws.Range("H15").Value = 38500
$ws.Range("I15").Value = 33000
$ws.Range("K15").Value = 33000
$ws.Range("M15").Value = -32712
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1319897.9
$ws.Range("I132").Value = 2120.1155
$ws.Range("K132").Value = 6360.3465
$ws.Range("M132").Value = -3830.3465
# Row 133 (Leve Item ID 41869)
$ws.Range("H133").Value = 79998.75
$ws.Range("J133").Value = 79998.75
$ws.Range("L133").Value = 79998.75
$ws.Range("N133").Value = -90118.75
